$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A2").Value = "ElementTD"
$ws.Range("B2").Value = "…"
$ws.Range("E2").Value = "Sommersemester 2017"
$ws.Range("G2").Value = "…"
$ws.Range("I2").Value = "GEBAUER Laurenz, PRANZ Bernhard, SCHILLER Markus"
$ws.Range("J2").Value = "MACHEINER Martin, BSc"
$ws.Range("K2").Value = "-"
$ws.Range("O2").Value = "…"
$ws.Range("S2").Value = "-"
